$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2-19, columns A-T (existing rows 2-13 get refreshed TPM values,
# and 6 new rows (14-19) are appended for Neutrophils as a new sending cluster).
$data = New-Object 'object[,]' 18,20
$data[0,0] = "FAPs"
$data[0,1] = "Csf2"
$data[0,2] = "Il3ra"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.2300786666666667
$data[0,7] = 0.690236
$data[0,8] = 0.2137022699341201
$data[0,9] = 0.2304051507653011
$data[0,10] = 2
$data[0,11] = 1
$data[0,12] = 5.1396195
$data[0,13] = 10.279239
$data[0,14] = 0.2307204726989479
$data[0,15] = 0.1721889818164019
$data[0,16] = 1.182516801734
$data[0,17] = 7.095100810404
$data[0,18] = 0.04930548873603836
$data[0,19] = 0.03967322831553177
$data[1,0] = "FAPs"
$data[1,1] = "Csf2"
$data[1,2] = "Il3ra"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.2300786666666667
$data[1,7] = 0.690236
$data[1,8] = 0.2137022699341201
$data[1,9] = 0.2304051507653011
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 4.164369
$data[1,13] = 12.493107
$data[1,14] = 0.1869409173525093
$data[1,15] = 0.2092737968300342
$data[1,16] = 0.9581324670279999
$data[1,17] = 8.623192203251998
$data[1,18] = 0.03994969838179797
$data[1,19] = 0.04821776070985101
$data[2,0] = "FAPs"
$data[2,1] = "Csf2"
$data[2,2] = "Il3ra"
$data[2,3] = "Inflammatory-Mac"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.2300786666666667
$data[2,7] = 0.690236
$data[2,8] = 0.2137022699341201
$data[2,9] = 0.2304051507653011
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 4.685485333333333
$data[2,13] = 14.056456
$data[2,14] = 0.2103341290013111
$data[2,15] = 0.2354616763543541
$data[2,16] = 1.078030218179556
$data[2,17] = 9.702271963615999
$data[2,18] = 0.04494888081219624
$data[2,19] = 0.05425158303987547
$data[3,0] = "FAPs"
$data[3,1] = "Csf2"
$data[3,2] = "Il3ra"
$data[3,3] = "MuSCs"
$data[3,4] = 2
$data[3,5] = 0.6666666666666666
$data[3,6] = 0.2300786666666667
$data[3,7] = 0.690236
$data[3,8] = 0.2137022699341201
$data[3,9] = 0.2304051507653011
$data[3,10] = 2
$data[3,11] = 1
$data[3,12] = 1.992128
$data[3,13] = 3.984256
$data[3,14] = 0.08942777064271192
$data[3,15] = 0.06674083401853875
$data[3,16] = 0.4583461540693333
$data[3,17] = 2.750076924416
$data[3,18] = 0.01911091758149541
$data[3,19] = 0.01537743192424336
$data[4,0] = "FAPs"
$data[4,1] = "Csf2"
$data[4,2] = "Il3ra"
$data[4,3] = "Neutrophils"
$data[4,4] = 2
$data[4,5] = 0.6666666666666666
$data[4,6] = 0.2300786666666667
$data[4,7] = 0.690236
$data[4,8] = 0.2137022699341201
$data[4,9] = 0.2304051507653011
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.585220333333333
$data[4,13] = 4.755661
$data[4,14] = 0.07116145166751166
$data[4,15] = 0.07966274793824443
$data[4,16] = 0.3647253806662222
$data[4,17] = 3.282528425996
$data[4,18] = 0.01520736375315442
$data[4,19] = 0.01835470744908939
$data[5,0] = "FAPs"
$data[5,1] = "Csf2"
$data[5,2] = "Il3ra"
$data[5,3] = "Resolving-Mac"
$data[5,4] = 2
$data[5,5] = 0.6666666666666666
$data[5,6] = 0.2300786666666667
$data[5,7] = 0.690236
$data[5,8] = 0.2137022699341201
$data[5,9] = 0.2304051507653011
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 4.709569
$data[5,13] = 14.128707
$data[5,14] = 0.2114152586370083
$data[5,15] = 0.2366719630424267
$data[5,16] = 1.083571356094667
$data[5,17] = 9.752142204852
$data[5,18] = 0.04517992066943778
$data[5,19] = 0.05453043932671008
$data[6,0] = "MuSCs"
$data[6,1] = "Csf2"
$data[6,2] = "Il3ra"
$data[6,3] = "ECs"
$data[6,4] = 1
$data[6,5] = 0.5
$data[6,6] = 0.2341465
$data[6,7] = 0.468293
$data[6,8] = 0.2174805655477089
$data[6,9] = 0.1563191709318771
$data[6,10] = 2
$data[6,11] = 1
$data[6,12] = 5.1396195
$data[6,13] = 10.279239
$data[6,14] = 0.2307204726989479
$data[6,15] = 0.1721889818164019
$data[6,16] = 1.20342391725675
$data[6,17] = 4.813695669027
$data[6,18] = 0.05017721888600192
$data[6,19] = 0.02691643888114402
$data[7,0] = "MuSCs"
$data[7,1] = "Csf2"
$data[7,2] = "Il3ra"
$data[7,3] = "FAPs"
$data[7,4] = 1
$data[7,5] = 0.5
$data[7,6] = 0.2341465
$data[7,7] = 0.468293
$data[7,8] = 0.2174805655477089
$data[7,9] = 0.1563191709318771
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 4.164369
$data[7,13] = 12.493107
$data[7,14] = 0.1869409173525093
$data[7,15] = 0.2092737968300342
$data[7,16] = 0.9750724260585
$data[7,17] = 5.850434556351
$data[7,18] = 0.04065601642983122
$data[7,19] = 0.03271350641823704
$data[8,0] = "MuSCs"
$data[8,1] = "Csf2"
$data[8,2] = "Il3ra"
$data[8,3] = "Inflammatory-Mac"
$data[8,4] = 1
$data[8,5] = 0.5
$data[8,6] = 0.2341465
$data[8,7] = 0.468293
$data[8,8] = 0.2174805655477089
$data[8,9] = 0.1563191709318771
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 4.685485333333333
$data[8,13] = 14.056456
$data[8,14] = 0.2103341290013111
$data[8,15] = 0.2354616763543541
$data[8,16] = 1.097089991601333
$data[8,17] = 6.582539949608001
$data[8,18] = 0.0457435853291899
$data[8,19] = 0.0368071740339426
$data[9,0] = "MuSCs"
$data[9,1] = "Csf2"
$data[9,2] = "Il3ra"
$data[9,3] = "MuSCs"
$data[9,4] = 1
$data[9,5] = 0.5
$data[9,6] = 0.2341465
$data[9,7] = 0.468293
$data[9,8] = 0.2174805655477089
$data[9,9] = 0.1563191709318771
$data[9,10] = 2
$data[9,11] = 1
$data[9,12] = 1.992128
$data[9,13] = 3.984256
$data[9,14] = 0.08942777064271192
$data[9,15] = 0.06674083401853875
$data[9,16] = 0.466449798752
$data[9,17] = 1.865799195008
$data[9,18] = 0.01944880213504779
$data[9,19] = 0.01043287184108
$data[10,0] = "MuSCs"
$data[10,1] = "Csf2"
$data[10,2] = "Il3ra"
$data[10,3] = "Neutrophils"
$data[10,4] = 1
$data[10,5] = 0.5
$data[10,6] = 0.2341465
$data[10,7] = 0.468293
$data[10,8] = 0.2174805655477089
$data[10,9] = 0.1563191709318771
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 1.585220333333333
$data[10,13] = 4.755661
$data[10,14] = 0.07116145166751166
$data[10,15] = 0.07966274793824443
$data[10,16] = 0.3711737927788333
$data[10,17] = 2.227042756673
$data[10,18] = 0.01547623275384639
$data[10,19] = 0.01245281471186147
$data[11,0] = "MuSCs"
$data[11,1] = "Csf2"
$data[11,2] = "Il3ra"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 1
$data[11,5] = 0.5
$data[11,6] = 0.2341465
$data[11,7] = 0.468293
$data[11,8] = 0.2174805655477089
$data[11,9] = 0.1563191709318771
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 4.709569
$data[11,13] = 14.128707
$data[11,14] = 0.2114152586370083
$data[11,15] = 0.2366719630424267
$data[11,16] = 1.1027290978585
$data[11,17] = 6.616374587151
$data[11,18] = 0.04597871001379172
$data[11,19] = 0.036996365045612
$data[12,0] = "Neutrophils"
$data[12,1] = "Csf2"
$data[12,2] = "Il3ra"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 0.6124066666666667
$data[12,7] = 1.83722
$data[12,8] = 0.5688171645181709
$data[12,9] = 0.6132756783028217
$data[12,10] = 2
$data[12,11] = 1
$data[12,12] = 5.1396195
$data[12,13] = 10.279239
$data[12,14] = 0.2307204726989479
$data[12,15] = 0.1721889818164019
$data[12,16] = 3.14753724593
$data[12,17] = 18.88522347558
$data[12,18] = 0.1312377650769076
$data[12,19] = 0.1055993146197261
$data[13,0] = "Neutrophils"
$data[13,1] = "Csf2"
$data[13,2] = "Il3ra"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 0.6124066666666667
$data[13,7] = 1.83722
$data[13,8] = 0.5688171645181709
$data[13,9] = 0.6132756783028217
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 4.164369
$data[13,13] = 12.493107
$data[13,14] = 0.1869409173525093
$data[13,15] = 0.2092737968300342
$data[13,16] = 2.55028733806
$data[13,17] = 22.95258604253999
$data[13,18] = 0.10633520254088
$data[13,19] = 0.1283425297019461
$data[14,0] = "Neutrophils"
$data[14,1] = "Csf2"
$data[14,2] = "Il3ra"
$data[14,3] = "Inflammatory-Mac"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 0.6124066666666667
$data[14,7] = 1.83722
$data[14,8] = 0.5688171645181709
$data[14,9] = 0.6132756783028217
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 4.685485333333333
$data[14,13] = 14.056456
$data[14,14] = 0.2103341290013111
$data[14,15] = 0.2354616763543541
$data[14,16] = 2.869422454702222
$data[14,17] = 25.82480209232
$data[14,18] = 0.119641662859925
$data[14,19] = 0.1444029192805359
$data[15,0] = "Neutrophils"
$data[15,1] = "Csf2"
$data[15,2] = "Il3ra"
$data[15,3] = "MuSCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 0.6124066666666667
$data[15,7] = 1.83722
$data[15,8] = 0.5688171645181709
$data[15,9] = 0.6132756783028217
$data[15,10] = 2
$data[15,11] = 1
$data[15,12] = 1.992128
$data[15,13] = 3.984256
$data[15,14] = 0.08942777064271192
$data[15,15] = 0.06674083401853875
$data[15,16] = 1.219992468053333
$data[15,17] = 7.319954808319999
$data[15,18] = 0.05086805092616872
$data[15,19] = 0.04093053025321539
$data[16,0] = "Neutrophils"
$data[16,1] = "Csf2"
$data[16,2] = "Il3ra"
$data[16,3] = "Neutrophils"
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 0.6124066666666667
$data[16,7] = 1.83722
$data[16,8] = 0.5688171645181709
$data[16,9] = 0.6132756783028217
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 1.585220333333333
$data[16,13] = 4.755661
$data[16,14] = 0.07116145166751166
$data[16,15] = 0.07966274793824443
$data[16,16] = 0.9707995002688888
$data[16,17] = 8.737195502419999
$data[16,18] = 0.04047785516051084
$data[16,19] = 0.04885522577729356
$data[17,0] = "Neutrophils"
$data[17,1] = "Csf2"
$data[17,2] = "Il3ra"
$data[17,3] = "Resolving-Mac"
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 0.6124066666666667
$data[17,7] = 1.83722
$data[17,8] = 0.5688171645181709
$data[17,9] = 0.6132756783028217
$data[17,10] = 3
$data[17,11] = 1
$data[17,12] = 4.709569
$data[17,13] = 14.128707
$data[17,14] = 0.2114152586370083
$data[17,15] = 0.2366719630424267
$data[17,16] = 2.884171452726667
$data[17,17] = 25.95754307454
$data[17,18] = 0.1202566279537788
$data[17,19] = 0.1451451586701046

$ws.Range("A2:T19").Value = $data

Write-Output "done"